$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.257.68'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '1.660.06'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.009'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.64%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '218.43'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5334'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('E7').Value = '  +0.59%  '
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06355'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.54'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07840'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.542'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.54%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.664.99'
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('D14').Value = '1.887.14'
$ws.Range('E14').Value = '  +0.48%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5517'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.13%  '
$ws.Range('D16').Value = '0.0₅8195'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('E17').Value = '  +0.79%  '
$ws.Range('D18').Value = '26.242.04'
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.647'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.13%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '192.11'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.78%  '
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('E23').Value = '  +0.97%  '
$ws.Range('E24').Value = '  +0.60%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '144.57'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.09%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1231'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.66%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.230'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '16.07'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.471'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.66%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05787'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.05%  '
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('E32').Value = '  +2.20%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.283'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.39%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.607'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.07%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.817'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.23%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.9574'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.56%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.429'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('E38').Value = '  +2.43%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01604'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.847'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('E41').Value = '  +1.10%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.008'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '104.68'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.94%  '
$ws.Range('D44').Value = '1.044.59'
$ws.Range('E44').Value = '  +3.76%  '
$ws.Range('D45').Value = '1.800.14'
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '57.10'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.42%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.012'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₈105'
$ws.Range('E48').Value = '  -1.80%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.4369'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.72%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.979'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.34%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05160'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.16%  '
